$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws = $wb.Worksheets.Item("Por Funcionalidad")

# --- Remove the "Por Tareas" sheet entirely -------------------------------
$wb.Worksheets.Item("Por Tareas").Delete() | Out-Null

# --- New "TOTAL" column (H) ------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 9.3

$ws.Range("H1").Value = "TOTAL"
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# --- New "TOTAL" row (23) --------------------------------------------------
$ws.Range("A23").Value = "TOTAL"
$ws.Range("A1").Copy()
$ws.Range("A23").PasteSpecial(-4122)

$ws.Range("B23:E23").Borders.LineStyle = 1

$ws.Range("F23").Formula = "=SUM(F2:F21)"
$ws.Range("G23").Formula = "=SUM(G2:G21)"
$ws.Range("H23").Formula = "=SUM(F23:G23)"

$excel.CutCopyMode = $false

# --- Selection matches the saved view in the edited workbook --------------
$ws.Activate()
$ws.Range("F29").Select() | Out-Null
